# Apply updated dSF (column F) values as part of a data repull / mean
# recalculation pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = 2
    10 = -1
    18 = -1
    24 = 5
    43 = 3
    46 = 0
    48 = 0
    50 = -2
    53 = -1
    54 = 2
    64 = -1
    65 = 0
    66 = 2
    70 = 0
    71 = -7
    72 = 3
    75 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
